$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update position ("D" column) picks that were re-drafted / corrected ---
$ws.Range("D199").Value = "SS"
$ws.Range("D440").Value = "MI"
$ws.Range("D441").Value = "DH"
$ws.Range("D462").Value = "DH"
$ws.Range("D468").Value = "3B"
$ws.Range("D470").Value = "DH"

# --- Append the final "end of bench" draft picks (rows 471-474) ---
# Copy formatting (currency style on C, date style on E) from the last
# existing data row so the new rows match the rest of the table.
$ws.Range("A470:E470").Copy()
$ws.Range("A471:E474").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A471").Value = "dembums"
$ws.Range("C471").Value = 1
$ws.Range("E471").Value = 43515

$ws.Range("A472").Value = "allrise"
$ws.Range("C472").Value = 1
$ws.Range("E472").Value = 43515

$ws.Range("A473").Value = "pasadena"
$ws.Range("C473").Value = 1
$ws.Range("E473").Value = 43515

$ws.Range("A474").Value = "allrise"
$ws.Range("C474").Value = 1
$ws.Range("E474").Value = 43515

# Player-name strings must land in the shared-string table in this order
# (Casey Mize, Mark Melancon , Victor Mesa, Josh Harrison) so they pick up
# the indices the workbook ends up with.
$ws.Range("B473").Value = "Casey Mize"
$ws.Range("D473").Value = "P"

$ws.Range("B474").Value = "Mark Melancon "
$ws.Range("D474").Value = "P"

$ws.Range("B472").Value = "Victor Mesa"
$ws.Range("D472").Value = "OF"

$ws.Range("B471").Value = "Josh Harrison"
$ws.Range("D471").Value = "2B"

# --- Restore the view state (scrolled position / active selection) ---
$ws.Range("D199").Select()
$excel.ActiveWindow.ScrollRow = 187
$excel.ActiveWindow.ScrollColumn = 1
